# Applies the "header content" update + new Kitchen data rows to the
# Excel Order form, per the commit:
#   "Excel Doc Complete (100%) - Change upon request on header content."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header block (rows 2-6): replace the "[Sample:] ..." placeholder
# values with real request data, and drop now-unused sample fields.
# ---------------------------------------------------------------------

# Row 2 - Contact / Data of Estimate / Sales Rep
$ws.Range("B2").Value = "DK NY"
$ws.Range("I2").ClearContents()
$ws.Range("M2").ClearContents()

# Row 3 - Firm / Project Number / Region
$ws.Range("B3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("M3").ClearContents()

# Row 4 - Address / Project Name
$ws.Range("B4").Value = "test@test.com"
$ws.Range("I4").ClearContents()

# Row 5 (new) - Email label, Project Address / Lead Time values
$ws.Range("A5").Value = "Email"
$ws.Range("A5").Font.FontStyle = "Bold"
$ws.Range("I5").Value = "affas"
$ws.Range("M5").Value = "8-10 WEEKS UPON ORDER APPROVAL AND RECEIPT OF DEPOSIT"

# Row 6 - Phone Number / Project Address value
$ws.Range("A6").Value = "Phone Number"
$ws.Range("B6").ClearContents()
$ws.Range("I6").Value = "afsfas"

# ---------------------------------------------------------------------
# New "Kitchen" line items (5 rows), matching the existing table
# layout used for the other items.
# ---------------------------------------------------------------------

$kitchenRows = 19..23
foreach ($r in $kitchenRows) {
    $ws.Cells.Item($r, 1).Value = "Kitchen"
    $ws.Cells.Item($r, 2).Value = "FE3008B-A1100010FF"
    $ws.Cells.Item($r, 3).Value = 1
    $ws.Cells.Item($r, 4).Value = "Canon de Fusil Anthracite"
    $ws.Cells.Item($r, 5).Value = "3.2"" x 3.2"" x .1"" (82 x 82 x 3 mm)"
    $ws.Cells.Item($r, 6).Value = "$"
    $ws.Cells.Item($r, 8).Value = "MELJAC"
    $ws.Cells.Item($r, 9).Value = "$"
    $ws.Cells.Item($r, 10).Value = 6
    $ws.Cells.Item($r, 11).Value = "$"
    $ws.Cells.Item($r, 12).Value = "$"
    $ws.Cells.Item($r, 13).Value = "$"
    $ws.Cells.Item($r, 14).Value = "$"
    $ws.Cells.Item($r, 15).Value = "TBD"
}

# Only the first row of the group carries the mechanism type/qty note
$ws.Cells.Item(19, 7).Value = "1 VV "
